$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values ---
$ws.Range("A2").Value = 150
$ws.Range("C2").Value = 100

# --- Normalize style of K1:L1 so they match the other headers (remove the ---
# --- now-unused fill style) by copying A1's formatting over them ---
$ws.Range("A1").Copy()
$ws.Range("K1:L1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Add new headers / data, writing cells in the exact order that makes ---
# --- the shared-string table grow as: OF_options, TAC, GHG, CAP, COP, OFs ---
$ws.Range("O1").Value = "OF_options"
$ws.Range("N2").Value = "TAC"
$ws.Range("N3").Value = "GHG"
$ws.Range("O3").Value = "CAP"
$ws.Range("P2").Value = "COP"
$ws.Range("N1").Value = "OFs"

$ws.Range("O2").Value = "TAC"
$ws.Range("P3").Value = "GHG"

# --- Update selection to match the saved view state ---
$ws.Range("O8").Select()
